# Update scripts with new TPM (transcripts-per-million) normalized values.
# Raw ligand/receptor expression for the "ECs" cluster was recomputed under
# the new TPM pipeline; every specificity / edge-weight column derived from
# it (I, J, O, P, Q, R, S, T) is refreshed accordingly across all 9 data rows.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 3.075165666666667
$ws.Range("H2").Value = 9.225497
$ws.Range("I2").Value = 0.02641273658732285
$ws.Range("J2").Value = 0.02641273658732285
$ws.Range("M2").Value = 0.989021
$ws.Range("N2").Value = 2.967063
$ws.Range("O2").Value = 0.05972921679266473
$ws.Range("P2").Value = 0.05972921679266473
$ws.Range("Q2").Value = 3.041403422812334
$ws.Range("R2").Value = 27.372630805311
$ws.Range("S2").Value = 0.001577612069711754
$ws.Range("T2").Value = 0.001577612069711754

$ws.Range("G3").Value = 3.075165666666667
$ws.Range("H3").Value = 9.225497
$ws.Range("I3").Value = 0.02641273658732285
$ws.Range("J3").Value = 0.02641273658732285
$ws.Range("O3").Value = 0.5654368392847325
$ws.Range("P3").Value = 0.5654368392847325
$ws.Range("Q3").Value = 28.79196531831922
$ws.Range("R3").Value = 259.127687864873
$ws.Range("S3").Value = 0.01493473429279604
$ws.Range("T3").Value = 0.01493473429279604

$ws.Range("G4").Value = 3.075165666666667
$ws.Range("H4").Value = 9.225497
$ws.Range("I4").Value = 0.02641273658732285
$ws.Range("J4").Value = 0.02641273658732285
$ws.Range("O4").Value = 0.3748339439226028
$ws.Range("P4").Value = 0.3748339439226028
$ws.Range("Q4").Value = 19.086492360845
$ws.Range("R4").Value = 171.778431247605
$ws.Range("S4").Value = 0.00990039022481505
$ws.Range("T4").Value = 0.00990039022481505

$ws.Range("I5").Value = 0.549422396165273
$ws.Range("J5").Value = 0.5494223961652731
$ws.Range("M5").Value = 0.989021
$ws.Range("N5").Value = 2.967063
$ws.Range("O5").Value = 0.05972921679266473
$ws.Range("P5").Value = 0.05972921679266473
$ws.Range("Q5").Value = 63.265506424989
$ws.Range("R5").Value = 569.389557824901
$ws.Range("S5").Value = 0.03281656941130091
$ws.Range("T5").Value = 0.03281656941130092

$ws.Range("I6").Value = 0.549422396165273
$ws.Range("J6").Value = 0.5494223961652731
$ws.Range("O6").Value = 0.5654368392847325
$ws.Range("P6").Value = 0.5654368392847325
$ws.Range("S6").Value = 0.3106636631199361
$ws.Range("T6").Value = 0.3106636631199362

$ws.Range("I7").Value = 0.549422396165273
$ws.Range("J7").Value = 0.5494223961652731
$ws.Range("O7").Value = 0.3748339439226028
$ws.Range("P7").Value = 0.3748339439226028
$ws.Range("S7").Value = 0.205942163634036
$ws.Range("T7").Value = 0.205942163634036

$ws.Range("I8").Value = 0.424164867247404
$ws.Range("J8").Value = 0.4241648672474041
$ws.Range("M8").Value = 0.989021
$ws.Range("N8").Value = 2.967063
$ws.Range("O8").Value = 0.05972921679266473
$ws.Range("P8").Value = 0.05972921679266473
$ws.Range("Q8").Value = 48.84221196913666
$ws.Range("R8").Value = 439.57990772223
$ws.Range("S8").Value = 0.02533503531165205
$ws.Range("T8").Value = 0.02533503531165205

$ws.Range("I9").Value = 0.424164867247404
$ws.Range("J9").Value = 0.4241648672474041
$ws.Range("O9").Value = 0.5654368392847325
$ws.Range("P9").Value = 0.5654368392847325
$ws.Range("S9").Value = 0.2398384418720003
$ws.Range("T9").Value = 0.2398384418720003

$ws.Range("I10").Value = 0.424164867247404
$ws.Range("J10").Value = 0.4241648672474041
$ws.Range("O10").Value = 0.3748339439226028
$ws.Range("P10").Value = 0.3748339439226028
$ws.Range("S10").Value = 0.1589913900637517
$ws.Range("T10").Value = 0.1589913900637517

Write-Output "applied"
